# Source data refresh for the EPEX Spot prices workbook.
#
# "Prix Spot" sheet: a new daily column ("16-dec") is inserted right before
# the existing "01-oct." column (at column EM), shifting every later column
# one to the right (old EM:FQ -> EN:FR). The freshly inserted column gets the
# header "16-dec" in row 1 and "-" placeholders for the 24 hourly data rows,
# matching the existing "no data yet" convention used elsewhere on the sheet.
#
# "Gaz" sheet: the last two days (2025-12-13 / 2025-12-14) get an updated
# price.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Prix Spot")
$ws.Columns("EM").Insert()
$ws.Range("EM1").Value = "16-dec"
$ws.Range("EM2:EM25").Value = "-"

$gaz = $wb.Worksheets.Item("Gaz")
$gaz.Range("B170").Value = 25.93
$gaz.Range("B171").Value = 25.93
